$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values, replacing the old "Strike#" values, for rows 2-31 (column G)
$kValues = @{
    2  = 0
    3  = 3
    4  = 4
    5  = 3
    6  = 9
    7  = 3
    8  = 8
    9  = 4
    10 = 8
    11 = 6
    12 = 8
    13 = 7
    14 = 6
    15 = 7
    16 = 5
    17 = 5
    18 = 7
    19 = 8
    20 = 4
    21 = 4
    22 = 7
    23 = 6
    24 = 4
    25 = 7
    26 = 5
    27 = 5
    28 = 5
    29 = 3
    30 = 3
    31 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
